$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 4216
$ws.Range("E2").Value = 273
$ws.Range("F2").Value = 273
$ws.Range("G2").Value = 263
$ws.Range("H2").Value = 150
$ws.Range("I2").Value = 150
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 7957
$ws.Range("L2").Value = 1654
$ws.Range("M2").Value = 6304
$ws.Range("N2").Value = 6303
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 275
$ws.Range("Q2").Value = 646
$ws.Range("R2").Value = -334
$ws.Range("S2").Value = -126
$ws.Range("T2").Value = 444
$ws.Range("U2").Value = 202
$ws.Range("V2").Value = 1111
$ws.Range("W2").Value = 6.48
$ws.Range("X2").Value = 3.55
$ws.Range("Y2").Value = 2.41
$ws.Range("Z2").Value = 1.88
$ws.Range("AA2").Value = 26.24
$ws.Range("AB2").Value = 2326.43
$ws.Range("AC2").Value = 273
$ws.Range("AD2").Value = 41.77
$ws.Range("AE2").Value = 11632
$ws.Range("AF2").Value = 0.98
$ws.Range("AG2").Value = 150
$ws.Range("AH2").Value = 1.32
$ws.Range("AI2").Value = 54.15
$ws.Range("AJ2").Value = 55000000

# Row 3
$ws.Range("D3").Value = 4071
$ws.Range("E3").Value = 353
$ws.Range("F3").Value = 353
$ws.Range("G3").Value = 247
$ws.Range("H3").Value = 125
$ws.Range("I3").Value = 125
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 7565
$ws.Range("L3").Value = 1165
$ws.Range("M3").Value = 6399
$ws.Range("N3").Value = 6398
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 275
$ws.Range("Q3").Value = 722
$ws.Range("R3").Value = -150
$ws.Range("S3").Value = -490
$ws.Range("T3").Value = 224
$ws.Range("U3").Value = 498
$ws.Range("V3").Value = 736
$ws.Range("W3").Value = 8.68
$ws.Range("X3").Value = 3.08
$ws.Range("Y3").Value = 1.97
$ws.Range("Z3").Value = 1.62
$ws.Range("AA3").Value = 18.21
$ws.Range("AB3").Value = 2342.56
$ws.Range("AC3").Value = 228
$ws.Range("AD3").Value = 54.24
$ws.Range("AE3").Value = 11808
$ws.Range("AF3").Value = 1.05
$ws.Range("AG3").Value = 200
$ws.Range("AH3").Value = 1.62
$ws.Range("AI3").Value = 86.54000000000001
$ws.Range("AJ3").Value = 55000000

# Row 4
$ws.Range("D4").Value = 4251
$ws.Range("E4").Value = 602
$ws.Range("F4").Value = 602
$ws.Range("G4").Value = 668
$ws.Range("H4").Value = 470
$ws.Range("I4").Value = 470
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 7509
$ws.Range("L4").Value = 789
$ws.Range("M4").Value = 6720
$ws.Range("N4").Value = 6716
$ws.Range("O4").Value = 4
$ws.Range("P4").Value = 275
$ws.Range("Q4").Value = 830
$ws.Range("R4").Value = -21
$ws.Range("S4").Value = -629
$ws.Range("T4").Value = 241
$ws.Range("U4").Value = 588
$ws.Range("V4").Value = 217
$ws.Range("W4").Value = 14.17
$ws.Range("X4").Value = 11.05
$ws.Range("Y4").Value = 7.16
$ws.Range("Z4").Value = 6.23
$ws.Range("AA4").Value = 11.74
$ws.Range("AB4").Value = 2470.65
$ws.Range("AC4").Value = 854
$ws.Range("AD4").Value = 16.98
$ws.Range("AE4").Value = 12394
$ws.Range("AF4").Value = 1.17
$ws.Range("AG4").Value = 500
$ws.Range("AH4").Value = 3.45
$ws.Range("AI4").Value = 57.68
$ws.Range("AJ4").Value = 55000000

# Row 5
$ws.Range("D5").Value = 4174
$ws.Range("E5").Value = 516
$ws.Range("F5").Value = 516
$ws.Range("G5").Value = 490
$ws.Range("H5").Value = 349
$ws.Range("I5").Value = 349
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 7247
$ws.Range("L5").Value = 858
$ws.Range("M5").Value = 6389
$ws.Range("N5").Value = 6385
$ws.Range("O5").Value = 4
$ws.Range("P5").Value = 275
$ws.Range("Q5").Value = 351
$ws.Range("R5").Value = -432
$ws.Range("S5").Value = -107
$ws.Range("T5").Value = 486
$ws.Range("U5").Value = -135
$ws.Range("V5").Value = 355
$ws.Range("W5").Value = 12.35
$ws.Range("X5").Value = 8.369999999999999
$ws.Range("Y5").Value = 5.33
$ws.Range("Z5").Value = 4.74
$ws.Range("AA5").Value = 13.42
$ws.Range("AB5").Value = 2497.97
$ws.Range("AC5").Value = 634
$ws.Range("AD5").Value = 42.95
$ws.Range("AE5").Value = 11784
$ws.Range("AF5").Value = 2.31
$ws.Range("AG5").Value = 130
$ws.Range("AH5").Value = 0.48
$ws.Range("AI5").Value = 20.19
$ws.Range("AJ5").Value = 55000000

# Row 6
$ws.Range("D6").Value = 4343
$ws.Range("E6").Value = 365
$ws.Range("F6").Value = 365
$ws.Range("G6").Value = 375
$ws.Range("H6").Value = 304
$ws.Range("I6").Value = 304
$ws.Range("K6").Value = 7394
$ws.Range("L6").Value = 788
$ws.Range("M6").Value = 6605
$ws.Range("N6").Value = 6604
$ws.Range("P6").Value = 275
$ws.Range("Q6").Value = 510
$ws.Range("R6").Value = -331
$ws.Range("S6").Value = -146
$ws.Range("T6").Value = 310
$ws.Range("U6").Value = 200
$ws.Range("V6").Value = 291
$ws.Range("W6").Value = 8.41
$ws.Range("X6").Value = 7.01
$ws.Range("Y6").Value = 4.68
$ws.Range("Z6").Value = 4.16
$ws.Range("AA6").Value = 11.93
$ws.Range("AB6").Value = 2512.63
$ws.Range("AC6").Value = 554
$ws.Range("AD6").Value = 37.39
$ws.Range("AE6").Value = 12188
$ws.Range("AF6").Value = 1.7
$ws.Range("AG6").Value = 80
$ws.Range("AH6").Value = 0.39
$ws.Range("AI6").Value = 14.26
$ws.Range("AJ6").Value = 54186600

# Row 7
$ws.Range("D7").Value = 4848
$ws.Range("E7").Value = 226
$ws.Range("G7").Value = 265
$ws.Range("H7").Value = 175
$ws.Range("I7").Value = 200
$ws.Range("K7").Value = 8165
$ws.Range("L7").Value = 1285
$ws.Range("M7").Value = 6880
$ws.Range("N7").Value = 6900
$ws.Range("P7").Value = 280
$ws.Range("Q7").Value = 330
$ws.Range("R7").Value = -390
$ws.Range("S7").Value = -210
$ws.Range("T7").Value = 300
$ws.Range("U7").Value = -35
$ws.Range("W7").Value = 4.67
$ws.Range("X7").Value = 3.61
$ws.Range("Y7").Value = 2.96
$ws.Range("Z7").Value = 2.25
$ws.Range("AA7").Value = 18.68
$ws.Range("AC7").Value = 369
$ws.Range("AD7").Value = 37.8
$ws.Range("AE7").Value = 12734
$ws.Range("AF7").Value = 1.1
$ws.Range("AG7").Value = 90
$ws.Range("AH7").Value = 0.65
$ws.Range("AI7").Value = 24.38

# Row 8
$ws.Range("D8").Value = 5418
$ws.Range("E8").Value = 468
$ws.Range("G8").Value = 485
$ws.Range("H8").Value = 355
$ws.Range("I8").Value = 379
$ws.Range("K8").Value = 8590
$ws.Range("L8").Value = 1400
$ws.Range("M8").Value = 7190
$ws.Range("N8").Value = 7230
$ws.Range("P8").Value = 280
$ws.Range("Q8").Value = 480
$ws.Range("R8").Value = -380
$ws.Range("S8").Value = -15
$ws.Range("T8").Value = 340
$ws.Range("U8").Value = 115
$ws.Range("W8").Value = 8.640000000000001
$ws.Range("X8").Value = 6.55
$ws.Range("Y8").Value = 5.36
$ws.Range("Z8").Value = 4.24
$ws.Range("AA8").Value = 19.47
$ws.Range("AC8").Value = 699
$ws.Range("AD8").Value = 19.94
$ws.Range("AE8").Value = 13343
$ws.Range("AF8").Value = 1.05
$ws.Range("AG8").Value = 100
$ws.Range("AH8").Value = 0.72
$ws.Range("AI8").Value = 14.3

# Row 9
$ws.Range("D9").Value = 5790
$ws.Range("E9").Value = 585
$ws.Range("G9").Value = 600
$ws.Range("H9").Value = 435
$ws.Range("I9").Value = 440
$ws.Range("K9").Value = 9085
$ws.Range("L9").Value = 1510
$ws.Range("M9").Value = 7575
$ws.Range("N9").Value = 7620
$ws.Range("P9").Value = 280
$ws.Range("Q9").Value = 720
$ws.Range("R9").Value = -395
$ws.Range("S9").Value = -10
$ws.Range("T9").Value = 355
$ws.Range("U9").Value = 340
$ws.Range("W9").Value = 10.1
$ws.Range("X9").Value = 7.51
$ws.Range("Y9").Value = 5.93
$ws.Range("Z9").Value = 4.92
$ws.Range("AA9").Value = 19.93
$ws.Range("AC9").Value = 812
$ws.Range("AD9").Value = 17.18
$ws.Range("AE9").Value = 14063
$ws.Range("AF9").Value = 0.99
$ws.Range("AG9").Value = 125
$ws.Range("AH9").Value = 0.9
$ws.Range("AI9").Value = 15.39

